$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3724909774500134
$ws.Range("C2").Value = 0.06096444224485253
$ws.Range("E2").Value = 0.573075119046436
$ws.Range("F2").Value = 2.493701241078199
$ws.Range("G2").Value = 0.7683712576383641
$ws.Range("H2").Value = 0.8580501211901179
$ws.Range("I2").Value = 0.6765997092239715
$ws.Range("J2").Value = 0.06361314572077958
$ws.Range("K2").Value = 0.4069831657943155

$ws.Range("B3").Value = 0.3329082712949116
$ws.Range("C3").Value = 0.05348938188573982
$ws.Range("E3").Value = 0.5471207742935036
$ws.Range("F3").Value = 2.464060770680263
$ws.Range("G3").Value = 0.7745014181456114
$ws.Range("H3").Value = 0.8661463243303231
$ws.Range("I3").Value = 0.6851358422237404
$ws.Range("J3").Value = 0.0642199734112765
$ws.Range("K3").Value = 0.3619623522568816

$ws.Range("B4").Value = 0.3086389125737981
$ws.Range("C4").Value = 0.04889948490702523
$ws.Range("E4").Value = 0.5314377078135806
$ws.Range("F4").Value = 2.447418340966095
$ws.Range("G4").Value = 0.7788880624012506
$ws.Range("H4").Value = 0.8715774164853158
$ws.Range("I4").Value = 0.6908180710564871
$ws.Range("J4").Value = 0.06463711258996341
$ws.Range("K4").Value = 0.3343497870479553

$ws.Range("B5").Value = 0.2987580555467275
$ws.Range("C5").Value = 0.0470290172477803
$ws.Range("E5").Value = 0.5251102939846675
$ws.Range("F5").Value = 2.441027092515341
$ws.Range("G5").Value = 0.7808317799306224
$ws.Range("H5").Value = 0.8739061917379658
$ws.Range("I5").Value = 0.6932443225027036
$ws.Range("J5").Value = 0.06481827404655505
$ws.Range("K5").Value = 0.323105447359211

$ws.Range("B6").Value = 0.2971179087263636
$ws.Range("C6").Value = 0.04671842462458642
$ws.Range("E6").Value = 0.524063469243174
$ws.Range("F6").Value = 2.439989402458963
$ws.Range("G6").Value = 0.7811639488418294
$ws.Range("H6").Value = 0.8742998595729219
$ws.Range("I6").Value = 0.6936538797380543
$ws.Range("J6").Value = 0.0648490298462967
$ws.Range("K6").Value = 0.3212388275771048

$ws.Range("B7").Value = 0.3085056184758059
$ws.Range("C7").Value = 0.04887425931921996
$ws.Range("E7").Value = 0.5313521168055644
$ws.Range("F7").Value = 2.447330565622153
$ws.Range("G7").Value = 0.7789136445260212
$ws.Range("H7").Value = 0.8716083554062592
$ws.Range("I7").Value = 0.6908503443895881
$ws.Range("J7").Value = 0.0646395105904265
$ws.Range("K7").Value = 0.3341981090556203

$ws.Range("B8").Value = 0.358835913338595
$ws.Range("C8").Value = 0.05838707486091721
$ws.Range("E8").Value = 0.5640735047185075
$ws.Range("F8").Value = 2.483157548012002
$ws.Range("G8").Value = 0.7703554186293644
$ws.Range("H8").Value = 0.8607461563766066
$ws.Range("I8").Value = 0.6794513177864907
$ws.Range("J8").Value = 0.06381311827521685
$ws.Range("K8").Value = 0.3914538817570303

$ws.Range("B9").Value = 0.4577951518533609
$ws.Range("C9").Value = 0.07704133622314657
$ws.Range("E9").Value = 0.6302541773728905
$ws.Range("F9").Value = 2.565812422018425
$ws.Range("G9").Value = 0.7585347379818899
$ws.Range("H9").Value = 0.8431001624634717
$ws.Range("I9").Value = 0.6606058955351308
$ws.Range("J9").Value = 0.06254725278392215
$ws.Range("K9").Value = 0.5039634792386494

$ws.Range("B10").Value = 0.5306506035614689
$ws.Range("C10").Value = 0.09074933775562499
$ws.Range("E10").Value = 0.6801200332765376
$ws.Range("F10").Value = 2.634169651478373
$ws.Range("G10").Value = 0.7529048902455884
$ws.Range("H10").Value = 0.8323706299216553
$ws.Range("I10").Value = 0.6489108220459556
$ws.Range("J10").Value = 0.06183520069965098
$ws.Range("K10").Value = 0.5867605303378696

$ws.Range("B11").Value = 0.5638255794814029
$ws.Range("C11").Value = 0.09698678098482105
$ws.Range("E11").Value = 0.7030788531841523
$ws.Range("F11").Value = 2.666940423114738
$ws.Range("G11").Value = 0.7510135272942193
$ws.Range("H11").Value = 0.8279763258801296
$ws.Range("I11").Value = 0.6440600335120727
$ws.Range("J11").Value = 0.06155898359595469
$ws.Range("K11").Value = 0.6244562695196123

$ws.Range("B12").Value = 0.5763924928410518
$ws.Range("C12").Value = 0.09934899860675728
$ws.Range("E12").Value = 0.7118124271806892
$ws.Range("F12").Value = 2.679591815170284
$ws.Range("G12").Value = 0.7503941364558671
$ws.Range("H12").Value = 0.8263824235449562
$ws.Range("I12").Value = 0.6422908823225413
$ws.Range("J12").Value = 0.06146127769405041
$ws.Range("K12").Value = 0.6387348975431735

$ws.Range("B13").Value = 0.5736858000371683
$ws.Range("C13").Value = 0.09884024235174138
$ws.Range("E13").Value = 0.7099297346516806
$ws.Range("F13").Value = 2.676856341027246
$ws.Range("G13").Value = 0.7505232185262543
$ws.Range("H13").Value = 0.8267225776808118
$ws.Range("I13").Value = 0.6426688841659143
$ws.Range("J13").Value = 0.06148201338722714
$ws.Range("K13").Value = 0.6356595645123662

$ws.Range("B14").Value = 0.5648593838584759
$ws.Range("C14").Value = 0.09718111735907087
$ws.Range("E14").Value = 0.7037965761057023
$ws.Range("F14").Value = 2.667976408784
$ws.Range("G14").Value = 0.7509606260499311
$ws.Range("H14").Value = 0.8278437882945013
$ws.Range("I14").Value = 0.6439131254263657
$ws.Range("J14").Value = 0.06155080696969151
$ws.Range("K14").Value = 0.6256309014349881

$ws.Range("B15").Value = 0.5594534925762389
$ws.Range("C15").Value = 0.09616488586053151
$ws.Range("E15").Value = 0.7000449947113765
$ws.Range("F15").Value = 2.662568714449208
$ws.Range("G15").Value = 0.7512411762272961
$ws.Range("H15").Value = 0.8285396990367673
$ws.Range("I15").Value = 0.6446840877511235
$ws.Range("J15").Value = 0.06159384347103014
$ws.Range("K15").Value = 0.6194885745379963

$ws.Range("B16").Value = 0.5284831615319661
$ws.Range("C16").Value = 0.09034173865066464
$ws.Range("E16").Value = 0.6786251527221339
$ws.Range("F16").Value = 2.632061788531388
$ws.Range("G16").Value = 0.7530420160879885
$ws.Range("H16").Value = 0.8326676128052526
$ws.Range("I16").Value = 0.6492372980222498
$ws.Range("J16").Value = 0.06185421473416142
$ws.Range("K16").Value = 0.5842976151973858

$ws.Range("B17").Value = 0.5094919272898721
$ws.Range("C17").Value = 0.08676983887778533
$ws.Range("E17").Value = 0.6655551592629507
$ws.Range("F17").Value = 2.613776415401361
$ws.Range("G17").Value = 0.7543186701441158
$ws.Range("H17").Value = 0.8353246946644788
$ws.Range("I17").Value = 0.6521509324926775
$ws.Range("J17").Value = 0.06202618397452042
$ws.Range("K17").Value = 0.5627167549792773

$ws.Range("B18").Value = 0.4985717714959321
$ws.Range("C18").Value = 0.08471552742042832
$ws.Range("E18").Value = 0.6580634854597918
$ws.Range("F18").Value = 2.60341672575808
$ws.Range("G18").Value = 0.7551159699306424
$ws.Range("H18").Value = 0.8368987712075153
$ws.Range("I18").Value = 0.653870941641415
$ws.Range("J18").Value = 0.06212958383631673
$ws.Range("K18").Value = 0.5503069500628612

$ws.Range("B19").Value = 0.4948749398843972
$ws.Range("C19").Value = 0.084019998701649
$ws.Range("E19").Value = 0.6555313692119853
$ws.Range("F19").Value = 2.599936146264554
$ws.Range("G19").Value = 0.7553967265550767
$ws.Range("H19").Value = 0.8374395869590501
$ws.Range("I19").Value = 0.6544608849209261
$ws.Range("J19").Value = 0.06216536303239195
$ws.Range("K19").Value = 0.5461057191311909

$ws.Range("B20").Value = 0.5115132588861115
$ws.Range("C20").Value = 0.08715005808255682
$ws.Range("E20").Value = 0.6669438073881366
$ws.Range("F20").Value = 2.615706609195144
$ws.Range("G20").Value = 0.7541762438310684
$ws.Range("H20").Value = 0.8350371032336739
$ws.Range("I20").Value = 0.6518361985876702
$ws.Range("J20").Value = 0.06200741286016154
$ws.Range("K20").Value = 0.5650137733934173

$ws.Range("B21").Value = 0.5674518037388054
$ws.Range("C21").Value = 0.09766843623151544
$ws.Range("E21").Value = 0.7055969589278277
$ws.Range("F21").Value = 2.670578089500879
$ws.Range("G21").Value = 0.750829516784151
$ws.Range("H21").Value = 0.8275125571731934
$ws.Range("I21").Value = 0.6435458212933085
$ws.Range("J21").Value = 0.06153041333386611
$ws.Range("K21").Value = 0.6285764554299362

$ws.Range("B22").Value = 0.6040355705868592
$ws.Range("C22").Value = 0.104544165790827
$ws.Range("E22").Value = 0.7310896968456149
$ws.Range("F22").Value = 2.707849705154587
$ws.Range("G22").Value = 0.7492068659627193
$ws.Range("H22").Value = 0.8230036455973249
$ws.Range("I22").Value = 0.6385225181088181
$ws.Range("J22").Value = 0.06125885079589111
$ws.Range("K22").Value = 0.6701418915353088

$ws.Range("B23").Value = 0.5845080060200019
$ws.Range("C23").Value = 0.1008743338620377
$ws.Range("E23").Value = 0.7174626085194262
$ws.Range("F23").Value = 2.687827804533839
$ws.Range("G23").Value = 0.7500210671546341
$ws.Range("H23").Value = 0.8253726799622427
$ws.Range("I23").Value = 0.6411673319841213
$ws.Range("J23").Value = 0.06140010120728334
$ws.Range("K23").Value = 0.6479556054248405

$ws.Range("B24").Value = 0.5105994205794389
$ws.Range("C24").Value = 0.08697816340819031
$ws.Range("E24").Value = 0.666315929683833
$ws.Range("F24").Value = 2.61483349263321
$ws.Range("G24").Value = 0.7542404375179075
$ws.Range("H24").Value = 0.8351669785258338
$ws.Range("I24").Value = 0.6519783498700096
$ws.Range("J24").Value = 0.06201588516585943
$ws.Range("K24").Value = 0.5639752996929133

$ws.Range("B25").Value = 0.4309970361667865
$ws.Range("C25").Value = 0.07199460237549715
$ws.Range("E25").Value = 0.6121333869644445
$ws.Range("F25").Value = 2.542117319470123
$ws.Range("G25").Value = 0.7611982482860782
$ws.Range("H25").Value = 0.8474818913949207
$ws.Range("I25").Value = 0.6653272775728212
$ws.Range("J25").Value = 0.06285155059450886
$ws.Range("K25").Value = 0.4735024081586232

Write-Output "Updated values for rows 2-25"